$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at row 1027, shifting existing rows 1027:1131 down to 1033:1137
$ws.Rows.Item(1027).Resize(6).Insert()

$ws.Cells.Item(1027,1).Value = 2
$ws.Cells.Item(1027,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(1027,3).Value = "Coquimbo"
$ws.Cells.Item(1027,4).Value = 44931
$ws.Cells.Item(1027,5).Value = 4
$ws.Cells.Item(1027,6).Value = 100112002
$ws.Cells.Item(1027,7).Value = "Pimiento"
$ws.Cells.Item(1027,8).Value = "Cuatro cascos rojo"
$ws.Cells.Item(1027,9).Value = "Primera"
$ws.Cells.Item(1027,10).Value = 340
$ws.Cells.Item(1027,11).Value = 7500
$ws.Cells.Item(1027,12).Value = 8000
$ws.Cells.Item(1027,13).Value = 7750
$ws.Cells.Item(1027,14).Value = "$/caja 18 kilos"
$ws.Cells.Item(1027,15).Value = "Provincia de Limarí"
$ws.Cells.Item(1027,16).Value = 431
$ws.Cells.Item(1027,17).Value = 18
$ws.Cells.Item(1027,18).Value = "Hortaliza"

$ws.Cells.Item(1028,1).Value = 2
$ws.Cells.Item(1028,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(1028,3).Value = "Coquimbo"
$ws.Cells.Item(1028,4).Value = 44931
$ws.Cells.Item(1028,5).Value = 4
$ws.Cells.Item(1028,6).Value = 100112002
$ws.Cells.Item(1028,7).Value = "Pimiento"
$ws.Cells.Item(1028,8).Value = "Cuatro cascos rojo"
$ws.Cells.Item(1028,9).Value = "Segunda"
$ws.Cells.Item(1028,10).Value = 400
$ws.Cells.Item(1028,11).Value = 5500
$ws.Cells.Item(1028,12).Value = 6000
$ws.Cells.Item(1028,13).Value = 5750
$ws.Cells.Item(1028,14).Value = "$/caja 18 kilos"
$ws.Cells.Item(1028,15).Value = "Provincia de Limarí"
$ws.Cells.Item(1028,16).Value = 319
$ws.Cells.Item(1028,17).Value = 18
$ws.Cells.Item(1028,18).Value = "Hortaliza"

$ws.Cells.Item(1029,1).Value = 2
$ws.Cells.Item(1029,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(1029,3).Value = "Coquimbo"
$ws.Cells.Item(1029,4).Value = 44931
$ws.Cells.Item(1029,5).Value = 4
$ws.Cells.Item(1029,6).Value = 100112002
$ws.Cells.Item(1029,7).Value = "Pimiento"
$ws.Cells.Item(1029,8).Value = "Cuatro cascos rojo"
$ws.Cells.Item(1029,9).Value = "Tercera"
$ws.Cells.Item(1029,10).Value = 400
$ws.Cells.Item(1029,11).Value = 3500
$ws.Cells.Item(1029,12).Value = 4000
$ws.Cells.Item(1029,13).Value = 3750
$ws.Cells.Item(1029,14).Value = "$/caja 18 kilos"
$ws.Cells.Item(1029,15).Value = "Provincia de Limarí"
$ws.Cells.Item(1029,16).Value = 208
$ws.Cells.Item(1029,17).Value = 18
$ws.Cells.Item(1029,18).Value = "Hortaliza"

$ws.Cells.Item(1030,1).Value = 2
$ws.Cells.Item(1030,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(1030,3).Value = "Coquimbo"
$ws.Cells.Item(1030,4).Value = 44931
$ws.Cells.Item(1030,5).Value = 4
$ws.Cells.Item(1030,6).Value = 100112002
$ws.Cells.Item(1030,7).Value = "Pimiento"
$ws.Cells.Item(1030,8).Value = "Cuatro cascos verde"
$ws.Cells.Item(1030,9).Value = "Primera"
$ws.Cells.Item(1030,10).Value = 700
$ws.Cells.Item(1030,11).Value = 6500
$ws.Cells.Item(1030,12).Value = 7000
$ws.Cells.Item(1030,13).Value = 6750
$ws.Cells.Item(1030,14).Value = "$/caja 18 kilos"
$ws.Cells.Item(1030,15).Value = "Provincia de Limarí"
$ws.Cells.Item(1030,16).Value = 375
$ws.Cells.Item(1030,17).Value = 18
$ws.Cells.Item(1030,18).Value = "Hortaliza"

$ws.Cells.Item(1031,1).Value = 2
$ws.Cells.Item(1031,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(1031,3).Value = "Coquimbo"
$ws.Cells.Item(1031,4).Value = 44931
$ws.Cells.Item(1031,5).Value = 4
$ws.Cells.Item(1031,6).Value = 100112002
$ws.Cells.Item(1031,7).Value = "Pimiento"
$ws.Cells.Item(1031,8).Value = "Cuatro cascos verde"
$ws.Cells.Item(1031,9).Value = "Segunda"
$ws.Cells.Item(1031,10).Value = 500
$ws.Cells.Item(1031,11).Value = 4500
$ws.Cells.Item(1031,12).Value = 5000
$ws.Cells.Item(1031,13).Value = 4750
$ws.Cells.Item(1031,14).Value = "$/caja 18 kilos"
$ws.Cells.Item(1031,15).Value = "Provincia de Limarí"
$ws.Cells.Item(1031,16).Value = 264
$ws.Cells.Item(1031,17).Value = 18
$ws.Cells.Item(1031,18).Value = "Hortaliza"

$ws.Cells.Item(1032,1).Value = 2
$ws.Cells.Item(1032,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(1032,3).Value = "Coquimbo"
$ws.Cells.Item(1032,4).Value = 44931
$ws.Cells.Item(1032,5).Value = 4
$ws.Cells.Item(1032,6).Value = 100112002
$ws.Cells.Item(1032,7).Value = "Pimiento"
$ws.Cells.Item(1032,8).Value = "Cuatro cascos verde"
$ws.Cells.Item(1032,9).Value = "Tercera"
$ws.Cells.Item(1032,10).Value = 200
$ws.Cells.Item(1032,11).Value = 2500
$ws.Cells.Item(1032,12).Value = 3000
$ws.Cells.Item(1032,13).Value = 2750
$ws.Cells.Item(1032,14).Value = "$/caja 18 kilos"
$ws.Cells.Item(1032,15).Value = "Provincia de Limarí"
$ws.Cells.Item(1032,16).Value = 153
$ws.Cells.Item(1032,17).Value = 18
$ws.Cells.Item(1032,18).Value = "Hortaliza"
